# Upload data from excel to db - append 10 new employee rows (109-118)
# following the same pattern as the existing rows (2-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$empIds = 109,110,111,112,113,114,115,116,117,118

$row = 10
foreach ($id in $empIds) {
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = "David"
    $ws.Cells.Item($row, 3).Value = 9578821821
    $ws.Cells.Item($row, 4).Value = 30
    $ws.Cells.Item($row, 6).Value = "Male"

    # Column G (DOJ) keeps the same date value/format as the rows above it -
    # copy the format from G2 so it reuses the existing date style instead
    # of creating a new one.
    $ws.Cells.Item($row, 7).Value = 40310
    $ws.Cells.Item(2, 7).Copy()
    $ws.Cells.Item($row, 7).PasteSpecial(-4122)

    $ws.Cells.Item($row, 8).Value = 50000

    $row = $row + 1
}
